$wb = $excel.ActiveWorkbook

$q = [char]34

# ---------------------------------------------------------------------------
# 1. InterpretatieCF: rename HPM.waarde(...) headers to HPM.waardeNum/waardeString
# ---------------------------------------------------------------------------
$wsCF = $wb.Worksheets.Item("InterpretatieCF")
$wsCF.Range("A1").Value = "HPM.waardeNum(" + $q + "IRT" + $q + ")"
$wsCF.Range("B1").Value = "HPM.waardeNum(" + $q + "PAP" + $q + ")"
$wsCF.Range("C1").Value = "HPM.waardeNum(" + $q + "DNA_N" + $q + ")"
$wsCF.Range("D1").Value = "HPM.waardeString(" + $q + "DNA_imut1" + $q + ")"
$wsCF.Range("E1").Value = "HPM.waardeString(" + $q + "DNA_imut2" + $q + ")"

# ---------------------------------------------------------------------------
# 2. InterpretatieEGABlad1: rename HPM.waarde(...) headers
# ---------------------------------------------------------------------------
$wsEGA1 = $wb.Worksheets.Item("InterpretatieEGABlad1")
$wsEGA1.Range("A1").Value = "HPM.waardeNum(" + $q + "EGA_N" + $q + ")"
$wsEGA1.Range("B1").Value = "HPM.waardeString(" + $q + "EGA_imut1" + $q + ")"
$wsEGA1.Range("C1").Value = "HPM.waardeString(" + $q + "EGA_imut2" + $q + ")"

# ---------------------------------------------------------------------------
# 3. InterpretatieEGABlad2: rename HPM.waarde(...) headers
# ---------------------------------------------------------------------------
$wsEGA2 = $wb.Worksheets.Item("InterpretatieEGABlad2")
$wsEGA2.Range("A1").Value = "HPM.waardeString(" + $q + "EGA_imut1" + $q + ")"
$wsEGA2.Range("B1").Value = "HPM.waardeString(" + $q + "EGA_imut2" + $q + ")"

# ---------------------------------------------------------------------------
# 4. Column widths: widen columns whose header text grew longer
# ---------------------------------------------------------------------------
$wsCF.Range("A1:E1").EntireColumn.AutoFit()
$wsEGA1.Range("A1:D1").ColumnWidth = 27.6666666666667
$wsEGA2.Range("A1:B1").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# 5. Sheet views / selections, and which sheet/tab is active
# ---------------------------------------------------------------------------
$wsEGA2.Activate()
$wsEGA2.Range("A1:C1048576").Select()

$wsEGA1.Activate()
$wsEGA1.Range("C1:C1048576").Select()

$wsCF.Activate()
$wsCF.Range("A1:F1048576").Select()
